# Test-data refresh: rename sheet/columns and insert a new "Unit Qty" column
# so the ERP export fixture looks like a real order-lines extract (adds qty
# support alongside the existing SKU/description/price columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Order Lines"

# Update header row (B1 -> Description, add Unit Qty column C, Unit Price moves to D)
$ws.Range("A1").Value = "Ordered Item"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Unit Qty"
$ws.Range("D1").Value = "Unit Price"

# Update data rows (SKU -> Ordered Item, Product Name -> Description, add Unit Qty, Unit Price shifts to D)
$ws.Range("A2").Value = "11080V012"
$ws.Range("B2").Value = "Premium 3-Piece Peeler Set"
$ws.Range("C2").Value = 24
$ws.Range("D2").Value = 12.99
$ws.Range("A3").Value = "11081V003"
$ws.Range("B3").Value = "Classic Swivel Peeler"
$ws.Range("C3").Value = 36
$ws.Range("D3").Value = 9.99
$ws.Range("A4").Value = "11082V008"
$ws.Range("B4").Value = "Ergonomic Can Opener"
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 21.99
$ws.Range("A5").Value = "11083V001"
$ws.Range("B5").Value = "Stainless Steel Garlic Press"
$ws.Range("C5").Value = 18
$ws.Range("D5").Value = 15.99
$ws.Range("A6").Value = "11084V005"
$ws.Range("B6").Value = "Locking Tongs 12-Inch"
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = 14.49
$ws.Range("A7").Value = "11085V002"
$ws.Range("B7").Value = "Silicone Spatula"
$ws.Range("C7").Value = 30
$ws.Range("D7").Value = 11.49
$ws.Range("A8").Value = "11086V010"
$ws.Range("B8").Value = "Balloon Whisk 11-Inch"
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 10.99
$ws.Range("A9").Value = "11087V004"
$ws.Range("B9").Value = "3-Piece Mixing Bowl Set"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 29.99
$ws.Range("A10").Value = "11088V006"
$ws.Range("B10").Value = "Salad Spinner Large"
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 32.99
$ws.Range("A11").Value = "11089V001"
$ws.Range("B11").Value = "Measuring Cups Set of 4"
$ws.Range("C11").Value = 36
$ws.Range("D11").Value = 9.49
$ws.Range("A12").Value = "11090V003"
$ws.Range("B12").Value = "Bamboo Cutting Board"
$ws.Range("C12").Value = 24
$ws.Range("D12").Value = 18.99
$ws.Range("A13").Value = "11091V007"
$ws.Range("B13").Value = "Stainless Steel Colander"
$ws.Range("C13").Value = 18
$ws.Range("D13").Value = 16.49
$ws.Range("A14").Value = "11092V002"
$ws.Range("B14").Value = "Box Grater 4-Sided"
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 13.99
$ws.Range("A15").Value = "11093V009"
$ws.Range("B15").Value = "Kitchen Shears Heavy Duty"
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 19.99
$ws.Range("A16").Value = "11094V001"
$ws.Range("B16").Value = "Flexible Turner"
$ws.Range("C16").Value = 24
$ws.Range("D16").Value = 12.49
$ws.Range("A17").Value = "11095V004"
$ws.Range("B17").Value = "Soup Ladle"
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = 11.99
$ws.Range("A18").Value = "11096V006"
$ws.Range("B18").Value = "Potato Masher"
$ws.Range("C18").Value = 18
$ws.Range("D18").Value = 13.49
$ws.Range("A19").Value = "11097V002"
$ws.Range("B19").Value = "Ice Cream Scoop"
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 10.99
$ws.Range("A20").Value = "11098V008"
$ws.Range("B20").Value = "Pizza Wheel Cutter"
$ws.Range("C20").Value = 24
$ws.Range("D20").Value = 14.99
$ws.Range("A21").Value = "11099V003"
$ws.Range("B21").Value = "Bottle Opener"
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 7.99

Write-Output "edit applied"
